$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values refreshed by the cryptos-list GitHub Actions run.
# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (e.g. "217.76", "2.400") are forced to stay plain text, matching the source data,
# by briefly switching the cell to the Text number format before assigning the value
# and then restoring the Normal style so no formatting changes leak into the sheet.

$ws.Range('D2').Value = '26.291.90'
$ws.Range('E2').Value = '  +1.05%  '
$ws.Range('D3').Value = '1.679.06'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5336'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.45%  '
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2681'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06475'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.94'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07542'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('D12').Value = '1.688.18'
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.525'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5777'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008453'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.79'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').Value = '26.317.52'
$ws.Range('E17').Value = '  +0.89%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.902'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.81%  '
$ws.Range('E20').Value = '  +0.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.25'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.207'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.008'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '145.69'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.40%  '
$ws.Range('E25').Value = '  +6.55%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.826'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.76'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06488'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.383'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.322'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.578'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.02%  '
$ws.Range('E32').Value = '  +1.68%  '
$ws.Range('E33').Value = '  +1.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.031'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6168'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.400'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.35%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.700'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.260'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('D39').Value = '1.111.38'
$ws.Range('E39').Value = '  +2.45%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01618'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.96%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8703'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.015'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.37'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').Value = '1.829.93'
$ws.Range('E44').Value = '  +0.70%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '57.07'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.33%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000107'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.159'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.32%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.001'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05266'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.57%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4289'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.076'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.13%  '
